$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldUrl = "https://github.com/shoot649854/IMG_DB/blob/main/profile.webp"
$newUrl = "https://raw.githubusercontent.com/shoot649854/IMG_DB/main/profile.webp"

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 17)  # Column Q = 17
    if ($cell.Value() -eq $oldUrl) {
        $cell.Value = $newUrl
    }
}
